# Daily update at 8 AM UTC
#
# Appends the next day's row of data to the sheet. The sheet keeps the
# last data row's date cell formatted differently (date+time format) from
# every other row's date cell (date-only format); when a new row is
# appended, that distinctive formatting moves down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the two number formats in play before making any changes:
#  - the format used by every "normal" (non-last) row in column A
#  - the distinct format currently applied to the last row (row 86)
$normalDateFormat = $ws.Range("A85").NumberFormat
$lastRowDateFormat = $ws.Range("A86").NumberFormat

# Row 86 is no longer the last row, so it switches to the normal format.
$ws.Range("A86").NumberFormat = $normalDateFormat

# Add the new row (87) with the next day's values.
$ws.Range("A87").Value = 45674
$ws.Range("B87").Value = 205
$ws.Range("C87").Value = 203
$ws.Range("D87").Value = 203

# Row 87 is now the last row, so it gets the distinct "last row" format.
$ws.Range("A87").NumberFormat = $lastRowDateFormat
